$wb = $excel.ActiveWorkbook
$wb.Windows.Item(1).ActiveSheet
